# Mark attendance as "A" (Absent) for the 9th Feb session for the
# students below. The target cells currently sit in column O (session 9)
# or column P (session 10) and are blank; they get the same style already
# used by existing "A" cells (style carries a text "A" shared-string plus
# the cell format Excel applies once a value is typed into it).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance Sheet")

# Reference cell that already holds the "Absent" mark/style so we can
# clone its formatting onto the newly-marked cells.
$absentTemplate = $ws.Range("L21")
$absentTemplate.Copy()

$cells = @(
    "P21",
    "P39",
    "O43",
    "O48",
    "P49",
    "O54",
    "O57",
    "P58",
    "O64",
    "P65",
    "O66",
    "O68",
    "P75",
    "O77",
    "P77",
    "O78",
    "O79",
    "O80"
)

foreach ($addr in $cells) {
    $target = $ws.Range($addr)
    $target.PasteSpecial(-4122)
    $target.Value = "A"
}

$excel.CutCopyMode = $false
$excel.Calculate()
